# Auto-generated Excel COM-interop script
# Applies the "Add data for 2024-07-07" update: refreshed 2024 (col K) totals
# and a couple of 2023 (col J) corrections across the citywide, by-neighborhood,
# and per-neighborhood violent-crime sheets.

$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item(1)
$ws.Range("J2").Value = 7710
$ws.Range("K2").Value = 4028
$ws.Range("K3").Value = 4126
$ws.Range("K4").Value = 829
$ws.Range("K5").Value = 293
$ws.Range("K6").Value = 4626
$ws.Range("J7").Value = 29290
$ws.Range("K7").Value = 13902

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item(2)
$ws.Range("K5").Value = 32
$ws.Range("K6").Value = 105
$ws.Range("K7").Value = 403
$ws.Range("K8").Value = 949
$ws.Range("K14").Value = 73
$ws.Range("K15").Value = 144
$ws.Range("K19").Value = 425
$ws.Range("K20").Value = 310
$ws.Range("K21").Value = 38
$ws.Range("K25").Value = 63
$ws.Range("K29").Value = 735
$ws.Range("K32").Value = 20
$ws.Range("K33").Value = 580
$ws.Range("K34").Value = 70
$ws.Range("K36").Value = 175
$ws.Range("K37").Value = 473
$ws.Range("K42").Value = 494
$ws.Range("K43").Value = 122
$ws.Range("K47").Value = 79
$ws.Range("K48").Value = 178
$ws.Range("K52").Value = 377
$ws.Range("K53").Value = 185
$ws.Range("K54").Value = 261
$ws.Range("K55").Value = 155
$ws.Range("K57").Value = 47
$ws.Range("J63").Value = 106
$ws.Range("K63").Value = 44
$ws.Range("K64").Value = 83
$ws.Range("K65").Value = 321
$ws.Range("K67").Value = 541
$ws.Range("K68").Value = 33
$ws.Range("K70").Value = 22
$ws.Range("K76").Value = 198
$ws.Range("K77").Value = 97
$ws.Range("K78").Value = 165
$ws.Range("K79").Value = 362
$ws.Range("K80").Value = 49
$ws.Range("K83").Value = 298
$ws.Range("K85").Value = 629
$ws.Range("K86").Value = 94
$ws.Range("K88").Value = 158
$ws.Range("K89").Value = 194
$ws.Range("K90").Value = 128
$ws.Range("K94").Value = 173
$ws.Range("K95").Value = 234
$ws.Range("K96").Value = 157
$ws.Range("K97").Value = 118
$ws.Range("K99").Value = 240
$ws.Range("K100").Value = 25
$ws.Range("J101").Value = 29290
$ws.Range("K101").Value = 13902

# Sheet 3: Bridgeport
$ws = $wb.Worksheets.Item(3)
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 73

# Sheet 4: West Ridge
$ws = $wb.Worksheets.Item(4)
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 157

# Sheet 5: Auburn Gresham
$ws = $wb.Worksheets.Item(5)
$ws.Range("K2").Value = 144
$ws.Range("K7").Value = 403

# Sheet 7: Uptown
$ws = $wb.Worksheets.Item(7)
$ws.Range("K2").Value = 52
$ws.Range("K7").Value = 194

# Sheet 8: South Shore
$ws = $wb.Worksheets.Item(8)
$ws.Range("K3").Value = 209
$ws.Range("K5").Value = 17
$ws.Range("K6").Value = 147
$ws.Range("K7").Value = 629

# Sheet 9: Little Village
$ws = $wb.Worksheets.Item(9)
$ws.Range("K4").Value = 22
$ws.Range("K6").Value = 148
$ws.Range("K7").Value = 377

# Sheet 11: Logan Square
$ws = $wb.Worksheets.Item(11)
$ws.Range("K3").Value = 42
$ws.Range("K7").Value = 185

# Sheet 12: Austin
$ws = $wb.Worksheets.Item(12)
$ws.Range("K2").Value = 267
$ws.Range("K3").Value = 282
$ws.Range("K4").Value = 54
$ws.Range("K6").Value = 320
$ws.Range("K7").Value = 949

# Sheet 13: South Chicago
$ws = $wb.Worksheets.Item(13)
$ws.Range("K2").Value = 105
$ws.Range("K3").Value = 103
$ws.Range("K7").Value = 298

# Sheet 14: Garfield Park
$ws = $wb.Worksheets.Item(14)
$ws.Range("K2").Value = 157
$ws.Range("K3").Value = 218
$ws.Range("K4").Value = 27
$ws.Range("K7").Value = 580

# Sheet 15: West Pullman
$ws = $wb.Worksheets.Item(15)
$ws.Range("K3").Value = 83
$ws.Range("K7").Value = 234

# Sheet 16: Grand Crossing
$ws = $wb.Worksheets.Item(16)
$ws.Range("K2").Value = 129
$ws.Range("K3").Value = 157
$ws.Range("K6").Value = 141
$ws.Range("K7").Value = 473

# Sheet 17: New City
$ws = $wb.Worksheets.Item(17)
$ws.Range("K3").Value = 81
$ws.Range("K7").Value = 321

# Sheet 18: Woodlawn
$ws = $wb.Worksheets.Item(18)
$ws.Range("K2").Value = 65
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 240

# Sheet 21: North Lawndale
$ws = $wb.Worksheets.Item(21)
$ws.Range("K2").Value = 158
$ws.Range("K3").Value = 188
$ws.Range("K6").Value = 159
$ws.Range("K7").Value = 541

# Sheet 24: Loop
$ws = $wb.Worksheets.Item(24)
$ws.Range("K6").Value = 127
$ws.Range("K7").Value = 261

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item(25)
$ws.Range("K2").Value = 209
$ws.Range("K3").Value = 261
$ws.Range("K4").Value = 38
$ws.Range("K6").Value = 206
$ws.Range("K7").Value = 735

# Sheet 26: Lake View
$ws = $wb.Worksheets.Item(26)
$ws.Range("K6").Value = 90
$ws.Range("K7").Value = 178

# Sheet 27: Chatham
$ws = $wb.Worksheets.Item(27)
$ws.Range("K2").Value = 137
$ws.Range("K3").Value = 128
$ws.Range("K4").Value = 17
$ws.Range("K7").Value = 425

# Sheet 29: River North
$ws = $wb.Worksheets.Item(29)
$ws.Range("K3").Value = 38
$ws.Range("K6").Value = 110
$ws.Range("K7").Value = 198

# Sheet 30: Ashburn
$ws = $wb.Worksheets.Item(30)
$ws.Range("K3").Value = 31
$ws.Range("K7").Value = 105

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item(32)
$ws.Range("K3").Value = 161
$ws.Range("K6").Value = 176
$ws.Range("K7").Value = 494

# Sheet 35: Rogers Park
$ws = $wb.Worksheets.Item(35)
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 165

# Sheet 36: Lower West Side
$ws = $wb.Worksheets.Item(36)
$ws.Range("K3").Value = 41
$ws.Range("K7").Value = 155

# Sheet 41: Chinatown
$ws = $wb.Worksheets.Item(41)
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 38

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item(42)
$ws.Range("K2").Value = 122
$ws.Range("K5").Value = 13
$ws.Range("K7").Value = 362

# Sheet 43: Near South Side
$ws = $wb.Worksheets.Item(43)
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 83

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item(44)
$ws.Range("K2").Value = 107
$ws.Range("K3").Value = 96
$ws.Range("K7").Value = 310

# Sheet 47: Grand Boulevard
$ws = $wb.Worksheets.Item(47)
$ws.Range("K2").Value = 68
$ws.Range("K3").Value = 50
$ws.Range("K7").Value = 175

# Sheet 49: Wrigleyville
$ws = $wb.Worksheets.Item(49)
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 25

# Sheet 50: Garfield Ridge
$ws = $wb.Worksheets.Item(50)
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 70

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item(51)
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 173

# Sheet 52: East Side
$ws = $wb.Worksheets.Item(52)
$ws.Range("K2").Value = 22
$ws.Range("K7").Value = 63

# Sheet 53: Kenwood
$ws = $wb.Worksheets.Item(53)
$ws.Range("K2").Value = 26
$ws.Range("K7").Value = 79

# Sheet 54: Brighton Park
$ws = $wb.Worksheets.Item(54)
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 144

# Sheet 65: West Town
$ws = $wb.Worksheets.Item(65)
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 118

# Sheet 67: O'Hare
$ws = $wb.Worksheets.Item(67)
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 22

# Sheet 68: United Center
$ws = $wb.Worksheets.Item(68)
$ws.Range("K3").Value = 47
$ws.Range("K7").Value = 158

# Sheet 69: Galewood
$ws = $wb.Worksheets.Item(69)
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 20

# Sheet 70: Armour Square
$ws = $wb.Worksheets.Item(70)
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 32

# Sheet 72: Streeterville
$ws = $wb.Worksheets.Item(72)
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 94

# Sheet 74: Washington Heights
$ws = $wb.Worksheets.Item(74)
$ws.Range("K3").Value = 41
$ws.Range("K7").Value = 128

# Sheet 76: North Park
$ws = $wb.Worksheets.Item(76)
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 33

# Sheet 77: Mckinley Park
$ws = $wb.Worksheets.Item(77)
$ws.Range("K4").Value = 3
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 47

# Sheet 79: Hyde Park
$ws = $wb.Worksheets.Item(79)
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 122

# Sheet 84: Riverdale
$ws = $wb.Worksheets.Item(84)
$ws.Range("K2").Value = 45
$ws.Range("K7").Value = 97

# Sheet 87: Rush & Division
$ws = $wb.Worksheets.Item(87)
$ws.Range("K2").Value = 13
$ws.Range("K7").Value = 49
